$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.50774764079014723
$ws.Range("B1").Value = 1.5447186534008406
$ws.Range("C1").Value = 0.77312631709436319
$ws.Range("D1").Value = 0.91732051115561342
$ws.Range("E1").Value = 1.0983302441261191
$ws.Range("F1").Value = -1.1029321864897761
$ws.Range("G1").Value = 1.753340040553288

$ws.Range("A2").Value = 0.3865631585471816
$ws.Range("B2").Value = 1.5309128263098699
$ws.Range("C2").Value = 0.74398068212453594
$ws.Range("D2").Value = 1.2072428800659996
$ws.Range("E2").Value = 1.4404079598246171
$ws.Range("F2").Value = -1.1919030721871431
$ws.Range("G2").Value = 2.1981944690401236

$ws.Range("A3").Value = 0.18867963690993386
$ws.Range("B3").Value = 1.4173982480063325
$ws.Range("C3").Value = 0.30372819600135698
$ws.Range("D3").Value = 1.3453011509757073
$ws.Range("E3").Value = 0.90965060721618518
$ws.Range("F3").Value = -1.1550875332778878
$ws.Range("G3").Value = 2.1245633912216131

$ws.Range("A4").Value = -0.059825250727540004
$ws.Range("B4").Value = 1.3882526130365052
$ws.Range("C4").Value = -0.19481556006147643
$ws.Range("D4").Value = 1.4266021327336462
$ws.Range("E4").Value = 0.047553404424454875
$ws.Range("F4").Value = -1.3284273623089653
$ws.Range("G4").Value = 1.8024274257656283

$ws.Range("A5").Value = -0.20095148321301901
$ws.Range("B5").Value = 1.4910293258248433
$ws.Range("C5").Value = -0.5905826033359719
$ws.Range("D5").Value = 1.2624661884298827
$ws.Range("E5").Value = -0.3175340230923277
$ws.Range("F5").Value = -1.3253594007331939
$ws.Range("G5").Value = 1.8315730607354554

$ws.Range("A6").Value = -0.69335931612430979
$ws.Range("B6").Value = 1.6106798272799232
$ws.Range("C6").Value = -0.61666027673002777
$ws.Range("D6").Value = 0.72250495109413704
$ws.Range("E6").Value = -0.86056322200384472
$ws.Range("F6").Value = -1.021631204731837
$ws.Range("G6").Value = 1.4342720366730746

$ws.Range("A7").Value = -0.75778650921550672
$ws.Range("B7").Value = 1.3836506706728484
$ws.Range("C7").Value = -0.53996123733574564
$ws.Range("D7").Value = 0.89891274170098567
$ws.Range("E7").Value = -1.0983302441261191
$ws.Range("F7").Value = -0.92345643430715596
$ws.Range("G7").Value = 1.227184630308513

$ws.Range("A8").Value = -0.56297094915403034
$ws.Range("B8").Value = 1.3222914391574228
$ws.Range("C8").Value = -0.52768939103266055
$ws.Range("D8").Value = 1.1658253987930873
$ws.Range("E8").Value = -1.1382137446111458
$ws.Range("F8").Value = -1.0446409165501216
$ws.Range("G8").Value = 1.1213399559444037

$ws.Range("A9").Value = -0.26844663787998718
$ws.Range("B9").Value = 1.1428156869748027
$ws.Range("C9").Value = -0.35281558121369749
$ws.Range("D9").Value = 1.5539225381281545
$ws.Range("E9").Value = -0.80227195206419033
$ws.Range("F9").Value = -1.1826991874598294
$ws.Range("G9").Value = 1.0400389741864646

$ws.Range("A10").Value = -0.22089323345553233
$ws.Range("B10").Value = 1.000155473701438
$ws.Range("C10").Value = 0.19634954084936207
$ws.Range("D10").Value = 1.7472041174017452
$ws.Range("E10").Value = -0.17794177139473438
$ws.Range("F10").Value = -1.1919030721871431
$ws.Range("G10").Value = 1.1259418983080607

$ws.Range("A11").Value = 0.14266021327336462
$ws.Range("B11").Value = 1.2716700731571966
$ws.Range("C11").Value = 0.6135923151542565
$ws.Range("D11").Value = 1.3529710549151355
$ws.Range("E11").Value = 0.82221370230670365
$ws.Range("F11").Value = -1.1919030721871431
$ws.Range("G11").Value = 1.6122138080678088

$ws.Range("A12").Value = 0.4878058905476339
$ws.Range("B12").Value = 1.3345632854605078
$ws.Range("C12").Value = 0.68415543139699597
$ws.Range("D12").Value = 0.94800012691332625
$ws.Range("E12").Value = 1.0523108204895499
$ws.Range("F12").Value = -0.89737876091310009
$ws.Range("G12").Value = 1.6398254622497503

$ws.Range("A13").Value = 0.56297094915403034
$ws.Range("B13").Value = 1.0906603401866908
$ws.Range("C13").Value = 0.51388356394168977
$ws.Range("D13").Value = 1.1244079175201749
$ws.Range("E13").Value = 1.1290098598838318
$ws.Range("F13").Value = -0.67034960430602519
$ws.Range("G13").Value = 1.7502720789775166

$ws.Range("A14").Value = 0.33287383097118411
$ws.Range("B14").Value = 0.95413605006486879
$ws.Range("C14").Value = 0.38196121618352463
$ws.Range("D14").Value = 1.4986992297642714
$ws.Range("E14").Value = 0.95413605006486879
$ws.Range("F14").Value = -0.81607777915516111
$ws.Range("G14").Value = 1.9665633700693921

$ws.Range("A15").Value = 0.11658253987930872
$ws.Range("B15").Value = 0.8958447801252144
$ws.Range("C15").Value = 0.17333982903107745
$ws.Range("D15").Value = 1.7073206169167185
$ws.Range("E15").Value = 0.59671852648751444
$ws.Range("F15").Value = -1.0630486860047492
$ws.Range("G15").Value = 1.8913983114629955

$ws.Range("A16").Value = -0.31293208072867079
$ws.Range("B16").Value = 0.91578653036772772
$ws.Range("C16").Value = -0.33133985018329848
$ws.Range("D16").Value = 1.5968740001889525
$ws.Range("E16").Value = -0.80994185600361857
$ws.Range("F16").Value = -0.98941760818623858
$ws.Range("G16").Value = 1.1121360712170898

$ws.Range("A17").Value = -0.60285444963905699
$ws.Range("B17").Value = 0.95720401164064006
$ws.Range("C17").Value = -0.25770877236478773
$ws.Range("D17").Value = 1.4312040750973032
$ws.Range("E17").Value = -0.94339818454966928
$ws.Range("F17").Value = -0.77312631709436319
$ws.Range("G17").Value = 0.99555353133778113

$ws.Range("A18").Value = -0.78079622103379132
$ws.Range("B18").Value = 1.1213399559444037
$ws.Range("C18").Value = -0.46786414030512058
$ws.Range("D18").Value = 1.021631204731837
$ws.Range("E18").Value = -1.1136700520049756
$ws.Range("F18").Value = -0.6688156235181395
$ws.Range("G18").Value = 1.2210487071569704

$ws.Range("A19").Value = -0.68722339297276724
$ws.Range("B19").Value = 0.94800012691332625
$ws.Range("C19").Value = -0.65347581563928314
$ws.Range("D19").Value = 1.0814564554593771
$ws.Range("E19").Value = -1.1566215140657734
$ws.Range("F19").Value = -0.69796125848796675
$ws.Range("G19").Value = 1.3453011509757073

$ws.Range("A20").Value = -0.4479223900626072
$ws.Range("B20").Value = 0.81300981757938984
$ws.Range("C20").Value = -0.58598066097231494
$ws.Range("D20").Value = 1.4588157292792447
$ws.Range("E20").Value = -1.0415729549743504
$ws.Range("F20").Value = -0.80073797127630464
$ws.Range("G20").Value = 1.1811652066719438

$ws.Range("A21").Value = -0.25924275315267337
$ws.Range("B21").Value = 0.65040785406351187
$ws.Range("C21").Value = -0.27458256103152978
$ws.Range("D21").Value = 1.7088545977046044
$ws.Range("E21").Value = -0.59365056491174317
$ws.Range("F21").Value = -0.84675739491287394
$ws.Range("G21").Value = 1.227184630308513

$ws.Range("A22").Value = -0.21322332951610412
$ws.Range("B22").Value = 0.56910687230557289
$ws.Range("C22").Value = 0.5905826033359719
$ws.Range("D22").Value = 1.7395342134623171
$ws.Range("E22").Value = 0.50774764079014723
$ws.Range("F22").Value = -0.69949523927585233
$ws.Range("G22").Value = 1.4787574795217582

$ws.Range("A23").Value = 0.24543692606170259
$ws.Range("B23").Value = 0.74091272054876467
$ws.Range("C23").Value = 0.51695152551746104
$ws.Range("D23").Value = 1.5569904997039259
$ws.Range("E23").Value = 0.94800012691332625
$ws.Range("F23").Value = -0.73477679739722213
$ws.Range("G23").Value = 1.727262367159232

$ws.Range("A24").Value = 0.5108156023659185
$ws.Range("B24").Value = 0.88664089539790059
$ws.Range("C24").Value = 0.5108156023659185
$ws.Range("D24").Value = 1.1919030721871431
$ws.Range("E24").Value = 1.047708878125893
$ws.Range("F24").Value = -0.55376706442671642
$ws.Range("G24").Value = 1.7303303287350031

$ws.Range("A25").Value = 0.40497092800180928
$ws.Range("B25").Value = 0.75849999999999995
$ws.Range("C25").Value = 0.7040971816395093
$ws.Range("D25").Value = 1.0952622825503477
$ws.Range("E25").Value = 1.0293011086712651
$ws.Range("F25").Value = -0.55547579448706197
$ws.Range("G25").Value = 1.385184651460734

$ws.Range("A26").Value = 0.27304858024364415
$ws.Range("B26").Value = 0.72613602891264795
$ws.Range("C26").Value = 0.65654377721505441
$ws.Range("D26").Value = 1.4312040750973032
$ws.Range("E26").Value = 1.1305438406717176
$ws.Range("F26").Value = -0.75690000000000002
$ws.Range("G26").Value = 1.701184693765176

$ws.Range("A27").Value = -0.08283496254582462
$ws.Range("B27").Value = 0.55069910285094514
$ws.Range("C27").Value = 0.27611654181941542
$ws.Range("D27").Value = 1.6587963479471199
$ws.Range("E27").Value = 0.31446606151655643
$ws.Range("F27").Value = -0.84675739491287394
$ws.Range("G27").Value = 1.5293788455219843

$ws.Range("A28").Value = -0.26231071472844464
$ws.Range("B28").Value = 0.54302919891151702
$ws.Range("C28").Value = -0.32673790781964157
$ws.Range("D28").Value = 1.6073788666742099
$ws.Range("E28").Value = -0.61972823830579904
$ws.Range("F28").Value = -0.73170883582145085
$ws.Range("G28").Value = 1.2778059963087391

$ws.Range("A29").Value = -0.40190296642603801
$ws.Range("B29").Value = 0.76085447079127799
$ws.Range("C29").Value = -0.59518454569962875
$ws.Range("D29").Value = 1.4388739790367313
$ws.Range("E29").Value = -0.9648739155800683
$ws.Range("F29").Value = -0.8390874909734457
$ws.Range("G29").Value = 1.3360972662483934

$ws.Range("A30").Value = -0.70869912400316626
$ws.Range("B30").Value = 0.88664089539790059
$ws.Range("C30").Value = -0.63200008460888413
$ws.Range("D30").Value = 1.0858836485505701
$ws.Range("E30").Value = -1.201106956914457
$ws.Range("F30").Value = -0.73477679739722213
$ws.Range("G30").Value = 1.3314953238847365

$ws.Range("D30").Select()
